$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "43.276.98"
$ws.Range("E2").Value = "  -1.14%  "
# Row 3
$ws.Range("D3").Value = "2.276.66"
$ws.Range("E3").Value = "  -0.55%  "
# Row 4
$ws.Range("E4").Value = "  +0.35%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.26"
$ws.Range("E5").Value = "  +0.42%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "264.67"
$ws.Range("E6").Value = "  -1.23%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.647"
$ws.Range("E7").Value = "  +3.57%  "
# Row 8
$ws.Range("E8").Value = "  +0.05%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("E9").Value = "  -1.16%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.52"
$ws.Range("E10").Value = "  -2.58%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0936"
$ws.Range("E11").Value = "  -1.65%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "9.20"
$ws.Range("E12").Value = "  +1.67%  "
# Row 13
$ws.Range("E13").Value = "  +1.11%  "
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "15.26"
$ws.Range("E14").Value = "  -3.14%  "
# Row 15
$ws.Range("D15").Value = "2.620.35"
$ws.Range("E15").Value = "  -0.45%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.860"
$ws.Range("E16").Value = "  +1.40%  "
# Row 17
$ws.Range("D17").Value = "2.276.59"
$ws.Range("E17").Value = "  -0.53%  "
# Row 18
$ws.Range("D18").Value = "43.246.37"
$ws.Range("E18").Value = "  -0.93%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0000110"
$ws.Range("E19").Value = "  +0.56%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.70"
$ws.Range("E20").Value = "  -0.55%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.31"
$ws.Range("E21").Value = "  +0.09%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.41"
$ws.Range("E22").Value = "  -1.54%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.10"
$ws.Range("E23").Value = "  +1.20%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.85"
$ws.Range("E24").Value = "  +3.29%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.33"
$ws.Range("E25").Value = "  -6.44%  "
# Row 26
$ws.Range("E26").Value = "  +1.80%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.52"
$ws.Range("E27").Value = "  -1.33%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "40.66"
$ws.Range("E28").Value = "  -2.87%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.35"
$ws.Range("E29").Value = "  -1.26%  "
# Row 30
$ws.Range("E30").Value = "  -1.75%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "173.29"
$ws.Range("E31").Value = "  -1.52%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.56"
$ws.Range("E32").Value = "  +0.22%  "
# Row 33
$ws.Range("E33").Value = "  -2.40%  "
# Row 34
$ws.Range("E34").Value = "  -0.04%  "
# Row 35
$ws.Range("E35").Value = "  +2.60%  "
# Row 36
$ws.Range("E36").Value = "  +4.12%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.70"
$ws.Range("E37").Value = "  -0.06%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.83"
$ws.Range("E38").Value = "  -0.62%  "
# Row 39
$ws.Range("E39").Value = "  -3.25%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.57"
$ws.Range("E40").Value = "  +7.42%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "14.65"
$ws.Range("E41").Value = "  +6.66%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "73.91"
$ws.Range("E42").Value = "  +0.76%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.234"
$ws.Range("E43").Value = "  -3.78%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.01"
$ws.Range("E44").Value = "  -4.88%  "
# Row 45
$ws.Range("E45").Value = "  -0.20%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.35"
$ws.Range("E46").Value = "  -1.52%  "
# Row 47
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.56"
$ws.Range("E47").Value = "  -2.38%  "
# Row 48
$ws.Range("B48").Value = "TrustWalletToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.26"
$ws.Range("E48").Value = "  +3.02%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0996"
$ws.Range("E49").Value = "  +0.67%  "
# Row 50
$ws.Range("B50").Value = "ordi"
$ws.Range("C50").Value = "https://coinranking.com/coin/j7-7vPrOi+ordi-ordi"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "71.52"
$ws.Range("E50").Value = "  +30.81%  "
# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "99.97"
$ws.Range("E51").Value = "  -2.09%  "
